# Updated cryptos list — refresh Price (D) and Volume(1h) (E) columns for
# the rows whose market data moved since the last snapshot.
#
# D and E cells hold plain text (not numbers), so for any D value that
# Excel's automatic type-detection would otherwise reinterpret as a number
# (losing formatting like trailing zeros, e.g. "1.010" -> 1.01, or collapsing
# a distinct cell type), we force the cell to Text before writing, then
# restore the default "Normal" style so no extra formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Cell,
        [string]$Value
    )
    $rng = $ws.Range($Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $Value
    $rng.Style = "Normal"
}

# row 2 - Bitcoin
Set-TextValue "D2" "30.535.79"
$ws.Range("E2").Value = "  +0.33%  "

# row 3 - Ethereum
Set-TextValue "D3" "2.113.23"
$ws.Range("E3").Value = "  +1.02%  "

# row 4 - TetherUSD
Set-TextValue "D4" "1.010"
$ws.Range("E4").Value = "  +0.70%  "

# row 5 - BNB
Set-TextValue "D5" "336.29"
$ws.Range("E5").Value = "  +1.96%  "

# row 7 - XRP
Set-TextValue "D7" "0.5240"
$ws.Range("E7").Value = "  +0.47%  "

# row 8 - Cardano
Set-TextValue "D8" "0.4549"
$ws.Range("E8").Value = "  +4.02%  "

# row 9 - OKB
Set-TextValue "D9" "54.94"
$ws.Range("E9").Value = "  +1.38%  "

# row 10 - Dogecoin
Set-TextValue "D10" "0.09113"
$ws.Range("E10").Value = "  +2.94%  "

# row 11 - Polygon (volume only)
$ws.Range("E11").Value = "  +1.51%  "

# row 12 - Solana
Set-TextValue "D12" "24.61"
$ws.Range("E12").Value = "  +1.69%  "

# row 13 - WrappedEther
Set-TextValue "D13" "2.113.16"
$ws.Range("E13").Value = "  +1.28%  "

# row 14 - Polkadot
Set-TextValue "D14" "6.844"
$ws.Range("E14").Value = "  +2.15%  "

# row 15 - Chainlink
Set-TextValue "D15" "8.120"
$ws.Range("E15").Value = "  +5.66%  "

# row 16 - ShibaInu
Set-TextValue "D16" "0.00001175"
$ws.Range("E16").Value = "  +4.86%  "

# row 17 - Litecoin
Set-TextValue "D17" "97.03"
$ws.Range("E17").Value = "  +1.25%  "

# row 18 - BinanceUSD (volume only)
$ws.Range("E18").Value = "  +0.63%  "

# row 19 - TRON
Set-TextValue "D19" "0.06680"
$ws.Range("E19").Value = "  +1.28%  "

# row 20 - Avalanche (volume only)
$ws.Range("E20").Value = "  +0.83%  "

# row 21 - Dai (volume only)
$ws.Range("E21").Value = "  +0.65%  "

# row 22 - Uniswap
Set-TextValue "D22" "6.280"
$ws.Range("E22").Value = "  +0.15%  "

# row 23 - WrappedBTC
Set-TextValue "D23" "30.613.90"
$ws.Range("E23").Value = "  +0.47%  "

# row 24 - Cosmos
Set-TextValue "D24" "12.80"
$ws.Range("E24").Value = "  +4.24%  "

# row 25 - Toncoin
Set-TextValue "D25" "2.356"
$ws.Range("E25").Value = "  +0.77%  "

# row 26 - Wrapped liquid staked Ether 2.0
Set-TextValue "D26" "2.362.03"
$ws.Range("E26").Value = "  +1.33%  "

# row 27 - EthereumClassic (volume only)
$ws.Range("E27").Value = "  +0.30%  "

# row 28 - Monero
Set-TextValue "D28" "163.86"
$ws.Range("E28").Value = "  +0.27%  "

# row 29 - LidoDAOToken
Set-TextValue "D29" "2.533"
$ws.Range("E29").Value = "  -1.17%  "

# row 30 - BitcoinCash
Set-TextValue "D30" "133.72"
$ws.Range("E30").Value = "  +1.64%  "

# row 31 - ImmutableX
Set-TextValue "D31" "1.213"
$ws.Range("E31").Value = "  +2.67%  "

# row 32 - Stellar
Set-TextValue "D32" "0.1070"
$ws.Range("E32").Value = "  +0.28%  "

# row 33 - ARBITRUM
Set-TextValue "D33" "1.639"
$ws.Range("E33").Value = "  -0.38%  "

# row 34 - Filecoin
Set-TextValue "D34" "6.363"
$ws.Range("E34").Value = "  +3.26%  "

# row 35 - HuobiToken (volume only)
$ws.Range("E35").Value = "  +1.13%  "

# row 36 - FraxShare (volume only)
$ws.Range("E36").Value = "  +4.77%  "

# row 37 - InternetComputer (DFINITY)
Set-TextValue "D37" "5.898"
$ws.Range("E37").Value = "  +8.02%  "

# row 38 - VeChain
Set-TextValue "D38" "0.02618"
$ws.Range("E38").Value = "  +1.70%  "

# row 39 - Hedera
Set-TextValue "D39" "0.06814"
$ws.Range("E39").Value = "  +0.11%  "

# row 40 - Algorand
Set-TextValue "D40" "0.2324"
$ws.Range("E40").Value = "  +3.09%  "

# row 41 - Aptos (volume only)
$ws.Range("E41").Value = "  -0.69%  "

# row 42 - TheSandbox (volume only)
$ws.Range("E42").Value = "  -0.09%  "

# row 43 - TrustWalletToken (volume only)
$ws.Range("E43").Value = "  +0.20%  "

# row 44 - EnergySwap
Set-TextValue "D44" "14.96"
$ws.Range("E44").Value = "  +7.63%  "

# row 45 - Decentraland
Set-TextValue "D45" "0.6442"
$ws.Range("E45").Value = "  +1.53%  "

# row 46 - NEARProtocol
Set-TextValue "D46" "2.309"
$ws.Range("E46").Value = "  +5.28%  "

# row 47 - PancakeSwap
Set-TextValue "D47" "3.688"
$ws.Range("E47").Value = "  +1.74%  "

# row 48 - BabyDogeCoin
Set-TextValue "D48" "0.00000000360"
$ws.Range("E48").Value = "  +21.58%  "

# row 49 - EOS
Set-TextValue "D49" "1.253"
$ws.Range("E49").Value = "  +0.71%  "

# row 50 - Aave
Set-TextValue "D50" "83.16"
$ws.Range("E50").Value = "  +1.68%  "

# row 51 - WOONetwork (volume only)
$ws.Range("E51").Value = "  +12.41%  "
